$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update justification text only
$ws.Range("F4").Value = "Osly Peralta"

# Row 5: update all fields (was "Permanente"/"Supervisor Legal"/... -> new values)
$ws.Range("A5").Value = "Temporal"
$ws.Range("B5").Value = "GESTOR DE COBROS"
$ws.Range("D5").Value = "Noviembre"
$ws.Range("E5").Value = "Diciembre"
$ws.Range("F5").Value = "Cierre de año 2021"

# C5 must remain text "2" (like the other quantity cells) rather than become a number.
# Force text formatting, set the value, then re-apply the original cell style (format only)
# from C4 so the style index matches the unchanged sibling cells.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2"
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6: remove entirely (dimension shrinks from A2:F6 to A2:F5)
$ws.Rows.Item(6).Delete()
